$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D data range as Text so numeric-looking price strings are not
# auto-converted to numbers when assigned below (matches original inline/shared string cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.420.88"
$ws.Range("E2").Value = "  +1.48%  "

$ws.Range("D3").Value = "1.864.51"
$ws.Range("E3").Value = "  +2.26%  "

$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "315.06"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("E6").Value = "  -0.30%  "

$ws.Range("D7").Value = "0.4639"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").Value = "0.3716"
$ws.Range("E8").Value = "  +1.53%  "

$ws.Range("D9").Value = "0.07352"
$ws.Range("E9").Value = "  +1.56%  "

$ws.Range("D10").Value = "0.8867"
$ws.Range("E10").Value = "  +3.01%  "

$ws.Range("D11").Value = "0.07929"
$ws.Range("E11").Value = "  +4.77%  "

$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").Value = "1.871.25"
$ws.Range("E13").Value = "  +2.93%  "

$ws.Range("D14").Value = "5.409"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "6.604"
$ws.Range("E15").Value = "  +1.39%  "

$ws.Range("D16").Value = "92.35"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").Value = "0.000008893"
$ws.Range("E18").Value = "  +2.93%  "

$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  -0.29%  "

$ws.Range("D20").Value = "14.87"
$ws.Range("E20").Value = "  +2.57%  "

$ws.Range("D21").Value = "27.456.70"
$ws.Range("E21").Value = "  +1.64%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  +0.25%  "

$ws.Range("D24").Value = "2.083.10"
$ws.Range("E24").Value = "  +4.79%  "

$ws.Range("D25").Value = "1.901"
$ws.Range("E25").Value = "  +3.11%  "

$ws.Range("D26").Value = "153.15"
$ws.Range("E26").Value = "  +1.03%  "

$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  +1.81%  "

$ws.Range("D28").Value = "2.075"
$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").Value = "5.139"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("D30").Value = "116.65"
$ws.Range("E30").Value = "  +1.21%  "

$ws.Range("D31").Value = "0.08888"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").Value = "3.028"
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").Value = "0.7548"
$ws.Range("E33").Value = "  +4.96%  "

$ws.Range("D34").Value = "1.165"
$ws.Range("E34").Value = "  +3.13%  "

$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("D36").Value = "2.656"
$ws.Range("E36").Value = "  +10.69%  "

$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "0.01962"
$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.999"
$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05252"
$ws.Range("E40").Value = "  -0.23%  "

$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").Value = "0.5185"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").Value = "0.1642"
$ws.Range("E43").Value = "  +0.89%  "

$ws.Range("D44").Value = "8.340"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").Value = "0.4855"
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("D46").Value = "10.34"
$ws.Range("E46").Value = "  +2.42%  "

$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").Value = "103.96"
$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").Value = "1.653"
$ws.Range("E49").Value = "  +1.70%  "

$ws.Range("D50").Value = "0.06250"
$ws.Range("E50").Value = "  +0.07%  "

$ws.Range("D51").Value = "65.75"
$ws.Range("E51").Value = "  +2.51%  "

# Restore default (no explicit number format) styling on column D so serialized cells
# remain style-free, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"